$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (column I) and IF (column J), matching the style of the
# existing header row (copy format from H1, which already has the header style).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for column I (I0) and column J (IF), one entry per data row
# (rows 2..72 of the sheet).
$i0Values = @(8,8,9,7,7,8,8,9,8,8,8,7,7,8,8,8,8,8,8,7,8,8,10,9,8,7,8,8,8,8,8,7,8,9,9,9,9,9,9,9,9,8,9,9,9,8,9,8,9,6,9,6,8,9,7,8,7,7,8,9,7,6,4,8,8,5,4,6,4,3,3)
$ifValues = @(9,8,9,7,8,8,8,9,8,8,8,8,7,8,8,8,8,8,8,7,8,8,10,9,8,8,8,8,8,8,8,8,8,9,9,9,9,9,9,9,9,9,9,9,9,8,9,8,9,7,9,7,8,9,7,8,8,7,8,9,7,7,4,8,8,6,5,6,4,3,3)

for ($k = 0; $k -lt $i0Values.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $i0Values[$k]
    $ws.Cells.Item($row, 10).Value = $ifValues[$k]
}
